$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Coby White'
$ws.Range("C2").Value = 'Chicago Bulls'

$ws.Range("A3").Value = 'Devin Booker'
$ws.Range("B3").Value = 'PG,SG'
$ws.Range("C3").Value = 'Phoenix Suns'

$ws.Range("A4").Value = 'Trae Young'
$ws.Range("C4").Value = 'Atlanta Hawks'

$ws.Range("A5").Value = 'Jamal Murray'
$ws.Range("B5").Value = 'PG,SG'
$ws.Range("C5").Value = 'Denver Nuggets'

$ws.Range("A6").Value = 'Luguentz Dort'
$ws.Range("C6").Value = 'Oklahoma City Thunder'

$ws.Range("A7").Value = 'P.J. Washington'
$ws.Range("B7").Value = 'PF'
$ws.Range("C7").Value = 'Dallas Mavericks'

$ws.Range("A8").Value = 'Clint Capela'
$ws.Range("B8").Value = 'C'
$ws.Range("C8").Value = 'Atlanta Hawks'

$ws.Range("A9").Value = 'Myles Turner'
$ws.Range("B9").Value = 'C'
$ws.Range("C9").Value = 'Indiana Pacers'

$ws.Range("A10").Value = 'LeBron James'
$ws.Range("B10").Value = 'SF,PF'
$ws.Range("C10").Value = 'Los Angeles Lakers'

$ws.Range("A11").Value = 'Kawhi Leonard'
$ws.Range("B11").Value = 'SG,SF,PF'
$ws.Range("C11").Value = 'LA Clippers'

$ws.Range("A12").Value = 'Devin Vassell'
$ws.Range("B12").Value = 'SG,SF'
$ws.Range("C12").Value = 'San Antonio Spurs'

$ws.Range("A13").Value = 'Jalen Brunson'
$ws.Range("B13").Value = 'PG'
$ws.Range("C13").Value = 'New York Knicks'

$ws.Range("A14").Value = 'Desmond Bane'
$ws.Range("B14").Value = 'SG,SF'
$ws.Range("C14").Value = 'Memphis Grizzlies'

$ws.Range("A15").Value = 'Norman Powell'
$ws.Range("B15").Value = 'SG,SF'
$ws.Range("C15").Value = 'LA Clippers'

$ws.Range("A16").Value = 'Walker Kessler'
$ws.Range("B16").Value = 'C'
$ws.Range("C16").Value = 'Utah Jazz'

$ws.Range("A19").Value = 'Immanuel Quickley'
$ws.Range("B19").Value = 'PG,SG'
$ws.Range("C19").Value = 'Toronto Raptors'
